$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 'I started working as a tutor.'
$ws.Range("G3").Value = 'Ah, so you got the job you applied for?'
$ws.Range("G5").Value = 'Ah, but don''t overdo it and push yourself too hard.'
$ws.Range("G6").Value = 'It''s fine. It''s only once a week.'
$ws.Range("G7").Value = 'Wow, that''s it?'
$ws.Range("G8").Value = 'Ah, she sounds surprised.'
$ws.Range("G10").Value = 'Plus, I can take days off whenever I want; there aren''t many part-time jobs this easygoing.'
$ws.Range("G11").Value = 'And I''m teaching an exam student.'
$ws.Range("G12").Value = 'Touya-kun… Are you really…?'
$ws.Range("G15").Value = 'Yuki, you don''t have to worry so much.'
$ws.Range("G16").Value = 'Even Yuki is telling me not to push myself too hard, huh.'
$ws.Range("G17").Value = 'It''s just that, Touya-kun… Sometimes you really do overdo it…'
$ws.Range("G18").Value = '…There are times when I wish I could be with you all the time…'
$ws.Range("G19").Value = 'It''s okay. I''m not overdoing it.'
$ws.Range("G20").Value = 'Okay…'
$ws.Range("G21").Value = 'I get it…'
$ws.Range("G22").Value = 'Whether it''s a lie or the truth,'
$ws.Range("G23").Value = 'there are times when it''s better to say things like this.'
$ws.Range("G24").Value = 'I promise. I won''t push myself to hard. I swear it.'
$ws.Range("G25").Value = 'Okay.'
$ws.Range("G26").Value = 'Do you feel relieved?'
$ws.Range("G27").Value = 'Yeah.'
$ws.Range("G28").Value = 'For example, there are times when just saying something like this is enough.'
$ws.Range("G29").Value = 'Even if it''s meaningless, there are times when words alone are desired.'
